$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "最近一次充电结束时间" (col D) timestamp for rows 2-50
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 4).Value = 45964.31144675926
}

# Update rows 18-50 with new station/terminal/time data
$ws.Cells.Item(18, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(18, 2).Value = '101号直流'
$ws.Cells.Item(18, 3).Value = 45954.028229166666
$ws.Cells.Item(19, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(19, 2).Value = '603号直流'
$ws.Cells.Item(19, 3).Value = 45959.03165509259
$ws.Cells.Item(20, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(20, 2).Value = '008B号直流'
$ws.Cells.Item(20, 3).Value = 45959.55945601852
$ws.Cells.Item(21, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(21, 2).Value = '702号直流'
$ws.Cells.Item(21, 3).Value = 45961.094305555554
$ws.Cells.Item(22, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(22, 2).Value = '103号直流'
$ws.Cells.Item(22, 3).Value = 45962.01840277778
$ws.Cells.Item(23, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(23, 2).Value = '002A号直流'
$ws.Cells.Item(23, 3).Value = 45962.15895833333
$ws.Cells.Item(24, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(24, 2).Value = '003B号直流'
$ws.Cells.Item(24, 3).Value = 45962.17201388889
$ws.Cells.Item(25, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(25, 2).Value = '705号直流'
$ws.Cells.Item(25, 3).Value = 45962.582824074074
$ws.Cells.Item(26, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(26, 2).Value = '904号直流'
$ws.Cells.Item(26, 3).Value = 45962.6437037037
$ws.Cells.Item(27, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(27, 2).Value = '102号直流'
$ws.Cells.Item(27, 3).Value = 45962.67511574074
$ws.Cells.Item(28, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(28, 2).Value = '703号直流'
$ws.Cells.Item(28, 3).Value = 45963.268900462965
$ws.Cells.Item(29, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(29, 2).Value = 'B01号直流'
$ws.Cells.Item(29, 3).Value = 45963.362488425926
$ws.Cells.Item(30, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(30, 2).Value = '901号直流'
$ws.Cells.Item(30, 3).Value = 45963.383564814816
$ws.Cells.Item(31, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(31, 2).Value = '206号直流'
$ws.Cells.Item(31, 3).Value = 45963.47337962963
$ws.Cells.Item(32, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(32, 2).Value = '602号直流'
$ws.Cells.Item(32, 3).Value = 45963.493946759256
$ws.Cells.Item(33, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(33, 2).Value = 'A03号直流'
$ws.Cells.Item(33, 3).Value = 45963.49891203704
$ws.Cells.Item(34, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(34, 2).Value = '109号直流'
$ws.Cells.Item(34, 3).Value = 45963.51157407407
$ws.Cells.Item(35, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(35, 2).Value = '004A号直流'
$ws.Cells.Item(35, 3).Value = 45963.52715277778
$ws.Cells.Item(36, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(36, 2).Value = '804号直流'
$ws.Cells.Item(36, 3).Value = 45963.528819444444
$ws.Cells.Item(37, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(37, 2).Value = '110号直流'
$ws.Cells.Item(37, 3).Value = 45963.5353125
$ws.Cells.Item(38, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(38, 2).Value = '205号直流'
$ws.Cells.Item(38, 3).Value = 45963.547002314815
$ws.Cells.Item(39, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(39, 2).Value = '107号直流'
$ws.Cells.Item(39, 3).Value = 45963.560428240744
$ws.Cells.Item(40, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(40, 2).Value = '201号直流'
$ws.Cells.Item(40, 3).Value = 45963.56891203704
$ws.Cells.Item(41, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(41, 2).Value = '402号直流'
$ws.Cells.Item(41, 3).Value = 45963.569710648146
$ws.Cells.Item(42, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(42, 2).Value = '203号直流'
$ws.Cells.Item(42, 3).Value = 45963.59641203703
$ws.Cells.Item(43, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(43, 2).Value = '104号直流'
$ws.Cells.Item(43, 3).Value = 45963.60113425926
$ws.Cells.Item(44, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(44, 2).Value = '302号直流'
$ws.Cells.Item(44, 3).Value = 45963.61158564815
$ws.Cells.Item(45, 1).Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Cells.Item(45, 2).Value = '108号直流'
$ws.Cells.Item(45, 3).Value = 45963.66019675926
$ws.Cells.Item(46, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(46, 2).Value = '905号直流'
$ws.Cells.Item(46, 3).Value = 45963.66409722222
$ws.Cells.Item(47, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(47, 2).Value = '504号直流'
$ws.Cells.Item(47, 3).Value = 45963.68914351852
$ws.Cells.Item(48, 1).Value = '长沙特来电飞狐四方坪西区充电站'
$ws.Cells.Item(48, 2).Value = '601号直流'
$ws.Cells.Item(48, 3).Value = 45963.72346064815
$ws.Cells.Item(49, 1).Value = '长沙特来电飞狐四方坪南区充电站'
$ws.Cells.Item(49, 2).Value = '105号直流'
$ws.Cells.Item(49, 3).Value = 45963.73883101852
$ws.Cells.Item(50, 1).Value = '长沙特来电飞狐四方坪东区充电站'
$ws.Cells.Item(50, 2).Value = '002B号直流'
$ws.Cells.Item(50, 3).Value = 45963.751539351855

# Clear rows 51-53 (no longer have data)
for ($r = 51; $r -le 53; $r++) {
    $ws.Cells.Item($r, 1).Value = $null
    $ws.Cells.Item($r, 2).Value = $null
    $ws.Cells.Item($r, 3).Value = $null
    $ws.Cells.Item($r, 4).Value = $null
}

# Update selection
$ws.Range("E8").Select()